$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# NOTE: the COM layer's ColumnWidth setter adds a constant ~0.8333333333333334
# (5/6) padding offset before it lands in the saved `width` attribute, so the
# assigned values below are pre-compensated (target - 0.8333333333333334) to
# land exactly on the target widths from the diff.
$ws.Columns.Item(4).ColumnWidth = 49.166666666666664
$ws.Columns.Item(6).ColumnWidth = 29.166666666666668
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 44.166666666666664
$ws.Columns.Item(9).ColumnWidth = 38.166666666666664
$ws.Columns.Item(10).ColumnWidth = 37.166666666666664
$ws.Columns.Item(11).ColumnWidth = 36.166666666666664
$ws.Columns.Item(13).ColumnWidth = 36.166666666666664

# --- Row 2: Student Government Association ---
$ws.Range("F2").Value = "studentgovernme@bladencc.edu"
$ws.Range("G2").Value = "(555) 890-1234"
$ws.Range("I2").Value = "https://instagram.com/studentgovernme"
$ws.Range("J2").Value = "https://facebook.com/studentgovernme"
$ws.Range("K2").Value = "https://twitter.com/studentgovernme"
$ws.Range("M2").Value = "https://tiktok.com/@studentgovernme"

# --- Row 3: Phi Theta Kappa Honor Society ---
$ws.Range("G3").Value = "(555) 123-4567"
$ws.Range("K3").Value = "https://twitter.com/phithetakappaho"

# --- Row 4: Student Volunteer Club ---
$ws.Range("F4").Value = "studentvoluntee@bladencc.edu"
$ws.Range("I4").Value = "https://instagram.com/studentvoluntee"
$ws.Range("J4").Value = "https://facebook.com/studentvoluntee"
$ws.Range("K4").Value = "https://twitter.com/studentvoluntee"

# --- Row 5: Future Teachers Association ---
$ws.Range("F5").Value = "futureteachersa@bladencc.edu"
$ws.Range("G5").Value = "(555) 123-4567"
$ws.Range("I5").Value = "https://instagram.com/futureteachersa"

# --- Row 6: Business Club ---
$ws.Range("J6").Value = "https://facebook.com/businessclub"

# --- Row 7: Art Club ---
$ws.Range("F7").Value = "artclub@bladencc.edu"
$ws.Range("G7").Value = "(555) 123-4567"
$ws.Range("J7").Value = "https://facebook.com/artclub"

# --- Row 8: Intramural Sports ---
$ws.Range("F8").Value = "intramuralsport@bladencc.edu"
$ws.Range("G8").Value = "(555) 789-0123"

# --- Row 9: International Student Association ---
$ws.Range("F9").Value = "internationalst@bladencc.edu"
$ws.Range("G9").Value = "(555) 789-0123"
$ws.Range("H9").Value = "https://linkedin.com/groups/internationalst"

# --- Row 10: Community Service Club ---
$ws.Range("D10").Value = "https://bladencc.edu/logos/communityservic_logo.png"
$ws.Range("G10").Value = "(555) 345-6789"
$ws.Range("I10").Value = "https://instagram.com/communityservic"
$ws.Range("K10").Value = "https://twitter.com/communityservic"
